$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting Late/heading/Outstanding
# one column to the right (N->O, O->P, P->Q), matching the "Loan RBI, Variable
# Instalments" layout change.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with a new selection, and move
# the selection on the "Transactions" sheet's active-tab flag off of it.
[void]$ws.Activate()
[void]$ws.Range("R8").Select()
